# retraining (student orgs 7)
#
# Add two new "patterns" rows to the student_orgs tag block on the
# "general" sheet (column B only, mirroring the existing pattern rows):
#   B44: "list all organizations in plm"
#   B45: "what various organizations can i join in plm?"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Write the new pattern values first.
$ws.Range("B44").Value = "list all organizations in plm"
$ws.Range("B45").Value = "what various organizations can i join in plm?"

# Copy the formatting (style) from the last existing pattern cell (B43)
# onto the two new cells so they match the rest of the column.
$fmtSrc = $ws.Range("B43")
$fmtSrc.Copy()
$ws.Range("B44:B45").PasteSpecial(-4122)
